# Weekly update: a new Damasco price entry for "Macroferia Regional de
# Talca" is inserted as the new row 13 (most recent week), pushing the
# previously-existing rows 13-17 down to rows 14-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13; existing rows 13:17 shift to 14:18.
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with this week's data.
$ws.Cells.Item(13, 1).Value = 5
$ws.Cells.Item(13, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(13, 3).Value = "Maule"
$ws.Cells.Item(13, 4).Value = 44529
$ws.Cells.Item(13, 5).Value = 7
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100103
$ws.Cells.Item(13, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(13, 9).Value = 100103003
$ws.Cells.Item(13, 10).Value = "Damasco"
$ws.Cells.Item(13, 11).Value = "Castle Brite"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 100
$ws.Cells.Item(13, 14).Value = 20000
$ws.Cells.Item(13, 15).Value = 20000
$ws.Cells.Item(13, 16).Value = 20000
$ws.Cells.Item(13, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(13, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(13, 19).Value = 1333
$ws.Cells.Item(13, 20).Value = 15

# Give the new date cell (D13) the same date format used by the other
# rows in this column ("D" cells use numFmtId 165: YYYY-MM-DD HH:MM:SS).
$ws.Cells.Item(13, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
